$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing H/I column values for rows 272-293 ---
$updates = @{
    272 = @{ H = 30644 }
    276 = @{ H = 13295 }
    277 = @{ H = 3222 }
    278 = @{ H = 29967 }
    279 = @{ H = 42294; I = 3080 }
    280 = @{ H = 35942 }
    281 = @{ H = 45614 }
    282 = @{ H = 46593; I = 2837 }
    284 = @{ H = 1095 }
    285 = @{ H = 40447 }
    286 = @{ H = 54756; I = 4187 }
    287 = @{ H = 57298; I = 3902 }
    288 = @{ H = 56222; I = 4102 }
    289 = @{ H = 64269; I = 4348 }
    290 = @{ H = 18292; I = 1501 }
    292 = @{ H = 78191; I = 6924 }
    293 = @{ H = 78922; I = 5663 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# --- Append new rows 294-298 ---
$newRows = @(
    @(44188, 161562, 114267, 45563, 14238, 2657, 1732, 86638, 4784),
    @(44189, 165608, 115663, 48213, 18443, 4046, 1732, 19236, 1136),
    @(44190, 166649, 116948, 47969, 4249, 1041, 1732, 1853, 145),
    @(44191, 167523, 119086, 46664, 3550, 874, 1773, 2568, 228),
    @(44192, 168092, 120410, 45803, 2540, 569, 1879, 2462, 250)
)

$startRow = 294
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
    # Column A carries the date style (numFmt yyyy-mm-dd) like the rest of the column.
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
}
